# Update cryptos list (Price + Volume(1h) columns) with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 48/49 swapped rank order (Aave now ranks above ordi) plus updated values.
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'

# Force Text format on Price cells before assigning, so numeric-looking strings
# (trailing zeros, thousands-dot groupings, etc.) are preserved verbatim instead
# of Excel auto-coercing them into numbers.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '41.208.90'
$ws.Range('E2').Value = '  -3.57%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.463.78'
$ws.Range('E3').Value = '  -2.56%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.87'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('E6').Value = '  -6.21%  '
$ws.Range('E7').Value = '  -2.79%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  -4.84%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '33.33'
$ws.Range('E10').Value = '  -6.11%  '
$ws.Range('E11').Value = '  -3.47%  '
$ws.Range('E12').Value = '  -1.13%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.98'
$ws.Range('E13').Value = '  -4.62%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.843.18'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '2.472.22'
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('E16').Value = '  -3.48%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.781'
$ws.Range('E17').Value = '  -3.84%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '41.164.68'
$ws.Range('E18').Value = '  -3.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '6.30'
$ws.Range('E19').Value = '  -5.62%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0921'
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.26'
$ws.Range('E21').Value = '  -8.15%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '68.31'
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '235.25'
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.75'
$ws.Range('E24').Value = '  -4.11%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('E26').Value = '  -6.07%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.99'
$ws.Range('E27').Value = '  -5.49%  '
$ws.Range('E28').Value = '  -6.01%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.61'
$ws.Range('E29').Value = '  -5.52%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '36.47'
$ws.Range('E30').Value = '  -5.25%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '152.08'
$ws.Range('E31').Value = '  -5.14%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.48'
$ws.Range('E32').Value = '  -5.08%  '
$ws.Range('E33').Value = '  -5.13%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.57'
$ws.Range('E34').Value = '  -4.23%  '
$ws.Range('E35').Value = '  -5.36%  '
$ws.Range('E36').Value = '  -2.49%  '
$ws.Range('E37').Value = '  -3.90%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '16.95'
$ws.Range('E38').Value = '  -7.87%  '
$ws.Range('E39').Value = '  -2.98%  '
$ws.Range('E40').Value = '  -7.79%  '
$ws.Range('E41').Value = '  +0.82%  '
$ws.Range('E42').Value = '  +0.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '19.98'
$ws.Range('E43').Value = '  -11.16%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.975.40'
$ws.Range('E44').Value = '  -1.26%  '
$ws.Range('E45').Value = '  -5.59%  '
$ws.Range('E46').Value = '  -8.01%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '8.67'
$ws.Range('E47').Value = '  -1.75%  '
$ws.Range('E50').Value = '  -6.58%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '74.05'
$ws.Range('E51').Value = '  -6.57%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "96.88"
$ws.Range("E48").Value = "  -3.75%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "68.90"
$ws.Range("E49").Value = "  -3.99%  "
